# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-09-18 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-19 Thursday", 2) | Out-Null

# Update the answer table. Several cells share identical source text
# (e.g. "60÷2=30, 0" appears twice with different targets), so each
# replacement is scoped to its own table cell rather than the whole
# document to avoid cross-contamination.
$t = $d.Tables.Item(1)

$edits = @(
    @{ Row = 1;  Col = 1; Old = "89÷9=9, 8";   New = "12÷9=1, 3" },
    @{ Row = 1;  Col = 2; Old = "60÷2=30, 0";  New = "10÷7=1, 3" },
    @{ Row = 1;  Col = 3; Old = "62÷7=8, 6";   New = "19÷6=3, 1" },
    @{ Row = 1;  Col = 4; Old = "60÷2=30, 0";  New = "70÷2=35, 0" },
    @{ Row = 1;  Col = 5; Old = "74÷8=9, 2";   New = "18÷5=3, 3" },

    @{ Row = 5;  Col = 1; Old = "56÷8=7, 0";   New = "40÷5=8, 0" },
    @{ Row = 5;  Col = 2; Old = "81÷8=10, 1";  New = "66÷9=7, 3" },
    @{ Row = 5;  Col = 3; Old = "87÷5=17, 2";  New = "65÷7=9, 2" },
    @{ Row = 5;  Col = 4; Old = "69÷2=34, 1";  New = "56÷3=18, 2" },
    @{ Row = 5;  Col = 5; Old = "37÷8=4, 5";   New = "95÷6=15, 5" },

    @{ Row = 9;  Col = 1; Old = "38÷3=12, 2";  New = "65÷9=7, 2" },
    @{ Row = 9;  Col = 2; Old = "15÷3=5, 0";   New = "70÷7=10, 0" },
    @{ Row = 9;  Col = 3; Old = "17÷9=1, 8";   New = "85÷7=12, 1" },
    @{ Row = 9;  Col = 4; Old = "88÷8=11, 0";  New = "22÷9=2, 4" },
    @{ Row = 9;  Col = 5; Old = "23÷7=3, 2";   New = "65÷9=7, 2" },

    @{ Row = 13; Col = 1; Old = "65÷5=13, 0";  New = "50÷7=7, 1" },
    @{ Row = 13; Col = 2; Old = "48÷8=6, 0";   New = "71÷9=7, 8" },
    @{ Row = 13; Col = 3; Old = "31÷7=4, 3";   New = "86÷9=9, 5" },
    @{ Row = 13; Col = 4; Old = "22÷8=2, 6";   New = "82÷4=20, 2" },
    @{ Row = 13; Col = 5; Old = "65÷3=21, 2";  New = "43÷4=10, 3" },

    @{ Row = 17; Col = 1; Old = "17÷2=8, 1";   New = "21÷7=3, 0" },
    @{ Row = 17; Col = 2; Old = "66÷8=8, 2";   New = "24÷8=3, 0" },
    @{ Row = 17; Col = 3; Old = "74÷6=12, 2";  New = "80÷6=13, 2" },
    @{ Row = 17; Col = 4; Old = "48÷7=6, 6";   New = "69÷9=7, 6" },
    @{ Row = 17; Col = 5; Old = "38÷7=5, 3";   New = "90÷6=15, 0" }
)

foreach ($e in $edits) {
    $cell = $t.Cell($e.Row, $e.Col)
    $cellRange = $cell.Range
    # Re-seat the cell's range through the document so mutations persist
    # (a bare $cell.Range is a disconnected snapshot for Find/replace).
    $scoped = $d.Range($cellRange.Start, $cellRange.End)
    $scoped.Find.Execute($e.Old, $true, $false, $false, $false, $false,
                          $true, 1, $false, $e.New, 2) | Out-Null
}
